$d = $word.ActiveDocument

# --- 1. "third year" -> "fourth year" --------------------------------------
$d.Content.Find.Execute("third year", $true, $false, $false, $false, $false,
                         $true, 1, $false, "fourth year", 2) | Out-Null

# --- 2. Relocate the "_GoBack" bookmark to the end of the "mentoring other
#        students." paragraph (right before the sentence we are about to
#        delete), THEN delete the trailing sentence. Doing the bookmark move
#        before the deletion keeps the Range anchors valid. --------------
$markerRange = $d.Content
$markerRange.Find.Execute("This September", $true, $false, $false, $false,
                           $false, $true, 1, $false, "", 0) | Out-Null
$markerRange.Collapse(1) | Out-Null

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null

$sentenceRange = $d.Content
$sentenceRange.Find.Execute("This September will be my second time partaking in a panel aimed at helping first year students transition into university.",
                             $true, $false, $false, $false, $false, $true, 1,
                             $false, "", 0) | Out-Null
$sentenceRange.Delete()

# --- 3. Append a trailing space to the final paragraph ---------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertAfter(" ")
